$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
}

$ws.Range("B3:B25").Select() | Out-Null
